$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select entire row 2 (as if the user clicked the row header) and clear its contents
$ws.Rows.Item(2).Select()
$ws.Range("A2:K2").ClearContents()
